$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 8320.5
$ws.Range("I28").Value = 1248.9
$ws.Range("J28").Value = 25999.5
$ws.Range("K28").Value = 1248.9
$ws.Range("L28").Value = 25999.5
$ws.Range("M28").Value = -763.9000000000001
$ws.Range("N28").Value = -26969.5
$ws.Range("H38").Value = 547.93335
$ws.Range("I38").Value = 170.6923
$ws.Range("K38").Value = 512.0769
$ws.Range("M38").Value = -140.0769
$ws.Range("H43").Value = 2090.1
$ws.Range("I43").Value = 1829.4
$ws.Range("J43").Value = 2350.8
$ws.Range("K43").Value = 1829.4
$ws.Range("L43").Value = 2350.8
$ws.Range("M43").Value = -1760.4
$ws.Range("N43").Value = -2488.8
$ws.Range("H70").Value = 3879.5
$ws.Range("I70").Value = 1031.6666
$ws.Range("K70").Value = 3094.9998
$ws.Range("M70").Value = -2824.9998
$ws.Range("H73").Value = 3879.5
$ws.Range("I73").Value = 1031.6666
$ws.Range("K73").Value = 3094.9998
$ws.Range("M73").Value = -2158.9998
$ws.Range("H125").Value = 2499.5
$ws.Range("I125").Value = 2499.5
$ws.Range("K125").Value = 22495.5
$ws.Range("M125").Value = -20035.5
$ws.Range("H132").Value = 15667.588
$ws.Range("I132").Value = 16988.846
$ws.Range("J132").Value = 11373.5
$ws.Range("K132").Value = 50966.538
$ws.Range("L132").Value = 34120.5
$ws.Range("M132").Value = -48436.538
$ws.Range("N132").Value = -39180.5
$ws.Range("H138").Value = 2026
$ws.Range("I138").Value = 1548
$ws.Range("J138").Value = 2663.3333
$ws.Range("K138").Value = 4644
$ws.Range("L138").Value = 7989.999899999999
$ws.Range("M138").Value = 496
$ws.Range("N138").Value = -18269.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 113
$ws.Range("I4").Value = 155.5
$ws.Range("K4").Value = 155.5
$ws.Range("M4").Value = -39.5
$ws.Range("H6").Value = 14000080
$ws.Range("I6").Value = 11666800
$ws.Range("J6").Value = 17500000
$ws.Range("K6").Value = 11666800
$ws.Range("L6").Value = 17500000
$ws.Range("M6").Value = -11666627
$ws.Range("N6").Value = -17500346
$ws.Range("H46").Value = 2625
$ws.Range("I46").Value = 3437
$ws.Range("J46").Value = 2393
$ws.Range("K46").Value = 3437
$ws.Range("L46").Value = 2393
$ws.Range("M46").Value = -3118
$ws.Range("N46").Value = -3031
$ws.Range("H110").Value = 4892.143
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("H132").Value = 2906.6667
$ws.Range("I132").Value = 2436.7083
$ws.Range("J132").Value = 6666.3335
$ws.Range("K132").Value = 7310.124899999999
$ws.Range("L132").Value = 19999.0005
$ws.Range("M132").Value = -4780.124899999999
$ws.Range("N132").Value = -25059.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 299999
$ws.Range("J27").Value = 299999
$ws.Range("L27").Value = 299999
$ws.Range("N27").Value = -300383
$ws.Range("H80").Value = 223.83333
$ws.Range("I80").Value = 184.6
$ws.Range("J80").Value = 238.92308
$ws.Range("K80").Value = 184.6
$ws.Range("L80").Value = 238.92308
$ws.Range("M80").Value = 813.4
$ws.Range("N80").Value = -2234.92308
$ws.Range("H83").Value = 223.83333
$ws.Range("I83").Value = 184.6
$ws.Range("J83").Value = 238.92308
$ws.Range("K83").Value = 923
$ws.Range("L83").Value = 1194.6154
$ws.Range("M83").Value = 4069
$ws.Range("N83").Value = -11178.6154
$ws.Range("H96").Value = 20666.334
$ws.Range("I96").Value = 20666.334
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 20666.334
$ws.Range("L96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -17920.334
$ws.Range("H99").Value = 2465.3333
$ws.Range("I99").Value = 2465.3333
$ws.Range("K99").Value = 2465.3333
$ws.Range("M99").Value = -967.3332999999998
$ws.Range("H107").Value = 3250.8276
$ws.Range("I107").Value = 1170.8889
$ws.Range("J107").Value = 6654.364
$ws.Range("K107").Value = 1170.8889
$ws.Range("L107").Value = 6654.364
$ws.Range("M107").Value = 749.1111000000001
$ws.Range("N107").Value = -10494.364

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 45000
$ws.Range("J20").Value = 45000
$ws.Range("L20").Value = 45000
$ws.Range("N20").Value = -45472
$ws.Range("H30").Value = 45000
$ws.Range("J30").Value = 45000
$ws.Range("L30").Value = 45000
$ws.Range("N30").Value = -45182
$ws.Range("H58").Value = 4691.0386
$ws.Range("J58").Value = 13348
$ws.Range("L58").Value = 13348
$ws.Range("N58").Value = -13754
$ws.Range("H123").Value = 39166
$ws.Range("J123").Value = 44999
$ws.Range("L123").Value = 44999
$ws.Range("N123").Value = -54799
$ws.Range("H128").Value = 45000
$ws.Range("J128").Value = 45000
$ws.Range("L128").Value = 45000
$ws.Range("N128").Value = -54960
$ws.Range("H134").Value = 2452.3044
$ws.Range("I134").Value = 2210.2727
$ws.Range("K134").Value = 6630.8181
$ws.Range("M134").Value = -4095.8181
$ws.Range("H136").Value = 4691.0386
$ws.Range("J136").Value = 13348
$ws.Range("L136").Value = 40044
$ws.Range("N136").Value = -45144

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 550
$ws.Range("I117").Value = 500
$ws.Range("J117").Value = 575
$ws.Range("K117").Value = 1500
$ws.Range("L117").Value = 1725
$ws.Range("M117").Value = 1942
$ws.Range("N117").Value = -8609
$ws.Range("H118").Value = 356.33334
$ws.Range("I118").Value = 356.33334
$ws.Range("K118").Value = 1069.00002
$ws.Range("M118").Value = 173.9999800000001
$ws.Range("H131").Value = 1000
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 22499.75
$ws.Range("J63").Value = 20000
$ws.Range("L63").Value = 20000
$ws.Range("N63").Value = -21372
$ws.Range("H66").Value = 22499.75
$ws.Range("J66").Value = 20000
$ws.Range("L66").Value = 60000
$ws.Range("N66").Value = -66864
$ws.Range("H70").Value = 4958.25
$ws.Range("I70").Value = 4380.857
$ws.Range("K70").Value = 4380.857
$ws.Range("M70").Value = -4110.857
$ws.Range("H73").Value = 4958.25
$ws.Range("I73").Value = 4380.857
$ws.Range("K73").Value = 4380.857
$ws.Range("M73").Value = -3444.857
$ws.Range("H113").Value = 7246.3335
$ws.Range("I113").Value = 5211.875
$ws.Range("K113").Value = 5211.875
$ws.Range("M113").Value = -3041.875
$ws.Range("H137").Value = 149759.8
$ws.Range("J137").Value = 174374.75
$ws.Range("L137").Value = 174374.75
$ws.Range("N137").Value = -184574.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6949.125
$ws.Range("I7").Value = 6299
$ws.Range("J7").Value = 8899.5
$ws.Range("K7").Value = 6299
$ws.Range("L7").Value = 8899.5
$ws.Range("M7").Value = -6187
$ws.Range("N7").Value = -9123.5
$ws.Range("H40").Value = 6857.0835
$ws.Range("I40").Value = 6329.1
$ws.Range("J40").Value = 9497
$ws.Range("K40").Value = 6329.1
$ws.Range("L40").Value = 9497
$ws.Range("M40").Value = -6193.1
$ws.Range("N40").Value = -9769
$ws.Range("H68").Value = 7583.3335
$ws.Range("J68").Value = 8699.799999999999
$ws.Range("L68").Value = 8699.799999999999
$ws.Range("N68").Value = -10197.8
$ws.Range("H71").Value = 7583.3335
$ws.Range("J71").Value = 8699.799999999999
$ws.Range("L71").Value = 43499
$ws.Range("N71").Value = -50987
$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996
$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984
$ws.Range("H126").Value = 6949.125
$ws.Range("I126").Value = 6299
$ws.Range("J126").Value = 8899.5
$ws.Range("K126").Value = 18897
$ws.Range("L126").Value = 26698.5
$ws.Range("M126").Value = -16427
$ws.Range("N126").Value = -31638.5
$ws.Range("H132").Value = 3462.0967
$ws.Range("I132").Value = 2470.0833
$ws.Range("K132").Value = 7410.249899999999
$ws.Range("M132").Value = -4880.249899999999
$ws.Range("H136").Value = 4932
$ws.Range("I136").Value = 4918.4
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 14755.2
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -12205.2
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 10000000
$ws.Range("J49").Value = 10000000
$ws.Range("L49").Value = 10000000
$ws.Range("N49").Value = -10000460
$ws.Range("H62").Value = 11000
$ws.Range("J62").Value = 11714.286
$ws.Range("L62").Value = 11714.286
$ws.Range("N62").Value = -12962.286
$ws.Range("H65").Value = 11000
$ws.Range("J65").Value = 11714.286
$ws.Range("L65").Value = 58571.43
$ws.Range("N65").Value = -64811.43
$ws.Range("H94").Value = 21150000
$ws.Range("J94").Value = 21150000
$ws.Range("L94").Value = 21150000
$ws.Range("N94").Value = -21151802
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H113").Value = 435
$ws.Range("I113").Value = 309.16666
$ws.Range("J113").Value = 812.5
$ws.Range("K113").Value = 927.4999799999999
$ws.Range("L113").Value = 2437.5
$ws.Range("M113").Value = 1242.50002
$ws.Range("N113").Value = -6777.5
$ws.Range("H136").Value = 5177.0454
$ws.Range("I136").Value = 3663
$ws.Range("J136").Value = 5744.8125
$ws.Range("K136").Value = 10989
$ws.Range("L136").Value = 17234.4375
$ws.Range("M136").Value = -8439
$ws.Range("N136").Value = -22334.4375
